$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New header column F: status_skl ---
$ws.Range("F1").Value = "status_skl"

# --- 2. Refresh the timestamp column (C) for the existing rows to today's date ---
$today = [DateTime]"2025-05-02"
$ws.Range("C2").Value = $today
$ws.Range("C3").Value = $today
$ws.Range("C2:C3").NumberFormat = "mm-dd-yy"

# --- 3. Replace the old formula-driven file_pdf column with plain text values ---
$ws.Range("E2").Value = "12345.pdf"
$ws.Range("E3").Value = "12346.pdf"

# --- 4. New "unduh skl" status column values for existing rows ---
$ws.Range("F2").Value = "LULUS"
$ws.Range("F3").Value = "DITAHAN"

# --- 5. Add the new 4th student row ---
$ws.Range("A4").Value = 12347
$ws.Range("B4").Value = "Lorem"
$ws.Range("C4").Value = $today
$ws.Range("C4").NumberFormat = "mm-dd-yy"
$ws.Range("D4").Value = "LULUS"
$ws.Range("E4").Value = "12347.pdf"
$ws.Range("F4").Value = "DITAHAN"

# --- 6. Box border (medium, light grey) around every used cell ---
$usedRange = $ws.Range("A1:F4")
$usedRange.Borders.Weight = -4138        # xlMedium
$usedRange.Borders.Color = 13421772      # RGB(204,204,204) == &HCCCCCC

$ws.Range("A1:F1").Font.Bold = $true
$ws.Range("A1:F1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A1:F1").VerticalAlignment = -4108     # xlCenter

$ws.Range("A1:F4").VerticalAlignment = -4108     # xlCenter
$ws.Range("F2:F4").HorizontalAlignment = -4108   # xlCenter (status_skl, like status_kelulusan)

# --- 7. Selection matches the new used range ---
$ws.Range("A1:F4").Select()
